$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# 1) Fix the "How to unlock" text on the "Ultra secret map" row (row 24, column C):
#    "Visit to different ultra secret rooms..." -> "Visit two different ultra secret rooms..."
$ws.Range("C24").Value = "Visit two different ultra secret rooms in a single run"

# 2) Highlight two additional rows with the green "reviewed" fill (same fill used on
#    rows 2, 14, 20, 47, 53-55, 60, 61): rows 10 and 13.
$greenColor = 5296274  # OLE/BGR encoding of RGB(146,208,80) = 0x92D050
$ws.Range("A10:C10").Interior.Color = $greenColor
$ws.Range("A13:C13").Interior.Color = $greenColor

# 3) Append a brand-new unlock entry as row 65 (previously a blank trailing row):
#    "Dad's empty wallet" | "Passive" | "Defeat a random modded enemy with Keeper"
$ws.Range("A65").Value = "Dad's empty wallet"
$ws.Range("B65").Value = "Passive"
$ws.Range("C65").Value = "Defeat a random modded enemy with Keeper"
$ws.Rows.Item(65).RowHeight = 41.4

# 4) Update the selection / scroll position to match where the author left off editing.
$ws.Range("C65").Select()
